# Generate Report for Handoff
#
# Two source files (828ddeb5-be3b-4b91-86db-c34f11b8562e.md and
# 82e2ea27-cdda-45e1-a5eb-c0c14a0dfb6c.md) have finished translation and are
# now ready to be handed off, but the existing handback file for each is
# stale relative to the latest source revision. Update the Overview sheet
# and the two per-locale sheets (zh-cn, de-de) to reflect the new status,
# refreshed "Latest Handoff" timestamps, and the stale-handback error detail.

$wb = $excel.ActiveWorkbook

$readyStatus = "Ready for handoff"

$errorFor828 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bdc31354d23f1453cd61641f3148ab4c50df16b8/e2e/828ddeb5-be3b-4b91-86db-c34f11b8562e.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/70cf525928629b42458fb4dc69264f51ab5f9876/e2e/828ddeb5-be3b-4b91-86db-c34f11b8562e.md."
$errorFor82e = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bdc31354d23f1453cd61641f3148ab4c50df16b8/e2e/82e2ea27-cdda-45e1-a5eb-c0c14a0dfb6c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/70cf525928629b42458fb4dc69264f51ab5f9876/e2e/82e2ea27-cdda-45e1-a5eb-c0c14a0dfb6c.md."

# --- Overview sheet: rows 4 & 5 are the two affected files ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E4").Value = $readyStatus
$overview.Range("F4").Value = $readyStatus
$overview.Range("G4").Value = "2016-08-22 16:27:32"

$overview.Range("E5").Value = $readyStatus
$overview.Range("F5").Value = $readyStatus
$overview.Range("G5").Value = "2016-08-22 16:27:32"

# --- zh-cn sheet: rows 4 (828ddeb5...) & 5 (82e2ea27...) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C4").Value = $readyStatus
$zhcn.Range("H4").Value = "2016-08-22 16:27:28"
$zhcn.Range("P4").Value = $errorFor828

$zhcn.Range("C5").Value = $readyStatus
$zhcn.Range("H5").Value = "2016-08-22 16:27:28"
$zhcn.Range("P5").Value = $errorFor82e

# --- de-de sheet: rows 4 (828ddeb5...) & 5 (82e2ea27...) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C4").Value = $readyStatus
$dede.Range("H4").Value = "2016-08-22 16:27:32"
$dede.Range("P4").Value = $errorFor828

$dede.Range("C5").Value = $readyStatus
$dede.Range("H5").Value = "2016-08-22 16:27:32"
$dede.Range("P5").Value = $errorFor82e

# Widen the Error Detail column (P) on both locale sheets so the new
# message is readable, matching the authored layout change.
# (39.2 "character" width units round-trips to the stored column width of
# exactly 40 in the underlying OOXML <col> element.)
$zhcn.Columns.Item(16).ColumnWidth = 39.2
$dede.Columns.Item(16).ColumnWidth = 39.2
